$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained a third data row (row 2 was the only match row; the
# scrape now also includes the Oct 23 2020 Sharjah match vs Mumbai Indians,
# which happens to repeat the same figures already present in row 2).
# All source cells are plain text (t="str"), including the numeric-looking
# ones (runs/balls/4s/6s/strike-rate), so force Text format first so Excel
# doesn't coerce "13", "10", "2", "0", "130.00" into real numbers.
$row3 = $ws.Range("A3:K3")
$row3.NumberFormat = "@"

$nbsp = [char]0x00A0

$ws.Range("A3").Value = " Sharjah"
$ws.Range("B3").Value = " October 23 2020"
$ws.Range("C3").Value = "Mumbai won by 10 wickets (with 46 balls remaining)"
$ws.Range("D3").Value = "Chennai Super Kings"
$ws.Range("E3").Value = "Mumbai Indians"
$ws.Range("F3").Value = "Imran Tahir" + $nbsp
$ws.Range("G3").Value = "13"
$ws.Range("H3").Value = "10"
$ws.Range("I3").Value = "2"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "130.00"
